$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

$ws.Range("B2").Value = 45208
$ws.Range("M2").Formula = "=ROUND((F2/`$D`$2-1)*100, 3)"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 45209
$ws.Range("C3").Value = 1
$ws.Range("D3").Formula = "=F2"
$ws.Range("E3").Value = 250
$ws.Range("F3").Formula = "=D3+E3"
$ws.Range("G3").Value = "ESPORTS"
$ws.Range("H3").Value = "WORLDS"
$ws.Range("I3").Value = "PSG"
$ws.Range("J3").Value = "GANA 1 MAPA EN LA SERIE"
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0
$ws.Range("M3").Formula = "=ROUND((F3/`$D`$2-1)*100, 3)"

$ws.Range("H7").Select()
